$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the two "<X>_ Petrelli2021_Cpx_Liq" typos (stray space after the
#    leading underscore) in the ThermoBar name column of the Cpx-Liq table.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("P_ Petrelli2021_Cpx_Liq", $true, $false, $false, `
    $false, $false, $true, 1, $false, "P_Petrelli2021_Cpx_Liq", 2) | Out-Null
$d.Content.Find.Execute("T_ Petrelli2021_Cpx_Liq", $true, $false, $false, `
    $false, $false, $true, 1, $false, "T_Petrelli2021_Cpx_Liq", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Each section heading ("Clinopyroxene-only Thermobarometers", "Two
#    pyroxene Thermobarometers", ...) is preceded by a run of visually
#    identical empty spacer paragraphs (same centered alignment & font
#    size as the heading). The edit drops one spacer paragraph from each
#    run and also strips the centered alignment from the (new) first
#    spacer paragraph of the run, leaving the rest of that paragraph's
#    formatting untouched.
# ---------------------------------------------------------------------------

function Find-ParaIndexByText($headingText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $headingText) {
            return $i
        }
    }
    return -1
}

function Remove-DuplicateSpacerParagraph($headingText) {
    $hidx = Find-ParaIndexByText($headingText)
    $hp = $d.Paragraphs.Item($hidx)
    $hsz = $hp.Range.Font.Size
    $halign = $hp.Alignment
    $hintable = $hp.Range.Information(12)

    $j = $hidx - 1
    $emptyStart = $hidx
    while ($j -ge 1) {
        $p = $d.Paragraphs.Item($j)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        $sz = $p.Range.Font.Size
        $align = $p.Alignment
        $intable = $p.Range.Information(12)
        if ($t.Length -eq 0 -and $sz -eq $hsz -and $align -eq $halign -and $intable -eq $hintable) {
            $emptyStart = $j
            $j = $j - 1
        } else {
            break
        }
    }

    # Drop the first spacer paragraph of the run entirely.
    $first = $d.Paragraphs.Item($emptyStart)
    $first.Range.Delete()

    # The paragraph that is now first in the (shortened) run loses its
    # centered alignment (its other paragraph/run formatting is kept).
    $newFirst = $d.Paragraphs.Item($emptyStart)
    $newFirst.Alignment = 0
}

# Process later headings first so earlier paragraph indices stay valid.
Remove-DuplicateSpacerParagraph("Two pyroxene Thermobarometers")
Remove-DuplicateSpacerParagraph("Clinopyroxene-only Thermobarometers")

# ---------------------------------------------------------------------------
# 3) Drop the trailing "Other Functions" header row and the
#    calculate_cpx_opx_press_temp(...) / calculate_cpx_opx_press_temp_matching(...)
#    description row from the last table.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Range.Text -like "*Other Functions*") {
        $rowCount = $t.Rows.Count
        $t.Rows.Item($rowCount).Delete()
        $t.Rows.Item($rowCount - 1).Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Update the trailing section's page size / margins.
# ---------------------------------------------------------------------------
$ps = $d.PageSetup
$ps.PageWidth = 595.45
$ps.PageHeight = 720.0
$ps.TopMargin = 12.95
$ps.BottomMargin = 12.95
$ps.LeftMargin = 12.95
$ps.RightMargin = 12.95
$ps.HeaderDistance = 35.3
$ps.FooterDistance = 35.3
$ps.Gutter = 0
